# Auto-generated: applies the Seraph_Profits market-data refresh.
# For each touched Leve row, updates currentAveragePrice / NQ / HQ columns
# (H-L) and the recomputed LeveProfitNQ/HQ columns (M, N) to match the
# latest scheduled-runner snapshot. Cells with no new profit value are
# cleared (matching the source row no longer carrying that column).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 3500
$ws.Range("J86").Value = 3500
$ws.Range("L86").Value = 3500
$ws.Range("N86").Value = -5746
$ws.Range("H89").Value = 3500
$ws.Range("J89").Value = 3500
$ws.Range("L89").Value = 17500
$ws.Range("N89").Value = -28732
$ws.Range("H112").Value = 2553.111
$ws.Range("J112").Value = 2712.5715
$ws.Range("L112").Value = 8137.7145
$ws.Range("N112").Value = -10353.7145
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 22332.166
$ws.Range("I46").Value = 23748.5
$ws.Range("K46").Value = 23748.5
$ws.Range("M46").Value = -23429.5
$ws.Range("H63").Value = 4044.111
$ws.Range("I63").Value = 3913.8572
$ws.Range("J63").Value = 4500
$ws.Range("K63").Value = 3913.8572
$ws.Range("L63").Value = 4500
$ws.Range("M63").Value = -3227.8572
$ws.Range("N63").Value = -5872
$ws.Range("H66").Value = 4044.111
$ws.Range("I66").Value = 3913.8572
$ws.Range("J66").Value = 4500
$ws.Range("K66").Value = 19569.286
$ws.Range("L66").Value = 22500
$ws.Range("M66").Value = -16137.286
$ws.Range("N66").Value = -29364
$ws.Range("H88").Value = 3000
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 3000
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 3000
$ws.Range("M88").Value = ""
$ws.Range("N88").Value = -3812
$ws.Range("H91").Value = 3000
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 3000
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 3000
$ws.Range("M91").Value = ""
$ws.Range("N91").Value = -5808
$ws.Range("H92").Value = 99500
$ws.Range("J92").Value = 99500
$ws.Range("L92").Value = 99500
$ws.Range("N92").Value = -104492
$ws.Range("H122").Value = 668941.9399999999
$ws.Range("I122").Value = 1112736.5
$ws.Range("K122").Value = 3338209.5
$ws.Range("M122").Value = -3335759.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4606
$ws.Range("I86").Value = 4358.75
$ws.Range("J86").Value = 4935.6665
$ws.Range("K86").Value = 4358.75
$ws.Range("L86").Value = 4935.6665
$ws.Range("M86").Value = -3235.75
$ws.Range("N86").Value = -7181.6665
$ws.Range("H89").Value = 4606
$ws.Range("I89").Value = 4358.75
$ws.Range("J89").Value = 4935.6665
$ws.Range("K89").Value = 21793.75
$ws.Range("L89").Value = 24678.3325
$ws.Range("M89").Value = -16177.75
$ws.Range("N89").Value = -35910.3325
$ws.Range("H99").Value = 2446.9678
$ws.Range("I99").Value = 2853.8667
$ws.Range("K99").Value = 2853.8667
$ws.Range("M99").Value = -1355.8667
$ws.Range("H134").Value = 2084.9
$ws.Range("I134").Value = 1776.4706
$ws.Range("K134").Value = 5329.4118
$ws.Range("M134").Value = -2794.4118
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 34789.77
$ws.Range("I62").Value = 4355.6665
$ws.Range("J62").Value = 399999
$ws.Range("K62").Value = 4355.6665
$ws.Range("L62").Value = 399999
$ws.Range("M62").Value = -3731.6665
$ws.Range("N62").Value = -401247
$ws.Range("H65").Value = 34789.77
$ws.Range("I65").Value = 4355.6665
$ws.Range("J65").Value = 399999
$ws.Range("K65").Value = 21778.3325
$ws.Range("L65").Value = 1999995
$ws.Range("M65").Value = -18658.3325
$ws.Range("N65").Value = -2006235
$ws.Range("H86").Value = 7817.6875
$ws.Range("I86").Value = 6989.9165
$ws.Range("K86").Value = 6989.9165
$ws.Range("M86").Value = -5866.9165
$ws.Range("H89").Value = 7817.6875
$ws.Range("I89").Value = 6989.9165
$ws.Range("K89").Value = 34949.5825
$ws.Range("M89").Value = -29333.5825
$ws.Range("H141").Value = 115850
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 12
$ws.Range("I6").Value = 13.4
$ws.Range("K6").Value = 40.2
$ws.Range("M6").Value = 72.8
$ws.Range("H12").Value = 737
$ws.Range("J12").Value = 736.1539
$ws.Range("L12").Value = 2208.4617
$ws.Range("N12").Value = -2554.4617
$ws.Range("H68").Value = 958.1429000000001
$ws.Range("I68").Value = 1062.5
$ws.Range("K68").Value = 3187.5
$ws.Range("M68").Value = -2376.5
$ws.Range("H71").Value = 958.1429000000001
$ws.Range("I71").Value = 1062.5
$ws.Range("K71").Value = 9562.5
$ws.Range("M71").Value = -5506.5
$ws.Range("H107").Value = 347.12
$ws.Range("I107").Value = 171.5
$ws.Range("K107").Value = 514.5
$ws.Range("M107").Value = 1405.5
$ws.Range("H117").Value = 2594.7144
$ws.Range("I117").Value = 900
$ws.Range("K117").Value = 2700
$ws.Range("M117").Value = 742
$ws.Range("H124").Value = 780
$ws.Range("I124").Value = 780
$ws.Range("K124").Value = 2340
$ws.Range("M124").Value = 2570
$ws.Range("H132").Value = 575
$ws.Range("I132").Value = 300
$ws.Range("J132").Value = 850
$ws.Range("K132").Value = 2700
$ws.Range("L132").Value = 7650
$ws.Range("M132").Value = -170
$ws.Range("N132").Value = -12710
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2648
$ws.Range("J7").Value = 2648
$ws.Range("L7").Value = 2648
$ws.Range("N7").Value = -2872
$ws.Range("H16").Value = 8599.200000000001
$ws.Range("I16").Value = 10499
$ws.Range("K16").Value = 10499
$ws.Range("M16").Value = -10329
$ws.Range("H36").Value = 80000
$ws.Range("J36").Value = 80000
$ws.Range("L36").Value = 80000
$ws.Range("N36").Value = -81124
$ws.Range("H46").Value = 3473.7856
$ws.Range("J46").Value = 4428.4287
$ws.Range("L46").Value = 4428.4287
$ws.Range("N46").Value = -4804.4287
$ws.Range("H82").Value = 2748.2222
$ws.Range("I82").Value = 3052
$ws.Range("J82").Value = 1958.4
$ws.Range("K82").Value = 3052
$ws.Range("L82").Value = 1958.4
$ws.Range("M82").Value = -2691
$ws.Range("N82").Value = -2680.4
$ws.Range("H85").Value = 2748.2222
$ws.Range("I85").Value = 3052
$ws.Range("J85").Value = 1958.4
$ws.Range("K85").Value = 3052
$ws.Range("L85").Value = 1958.4
$ws.Range("M85").Value = -1804
$ws.Range("N85").Value = -4454.4
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").Value = ""
$ws.Range("H100").Value = 1876.3077
$ws.Range("I100").Value = 1710.2222
$ws.Range("J100").Value = 2250
$ws.Range("K100").Value = 1710.2222
$ws.Range("L100").Value = 2250
$ws.Range("M100").Value = -1169.2222
$ws.Range("N100").Value = -3332
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").Value = ""
$ws.Range("H122").Value = 5951.4614
$ws.Range("I122").Value = 5874.25
$ws.Range("K122").Value = 17622.75
$ws.Range("M122").Value = -15172.75
$ws.Range("H126").Value = 2648
$ws.Range("J126").Value = 2648
$ws.Range("L126").Value = 7944
$ws.Range("N126").Value = -12884
$ws.Range("H133").Value = 125000
$ws.Range("J133").Value = 125000
$ws.Range("L133").Value = 125000
$ws.Range("N133").Value = -130060
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H115").Value = 40000
$ws.Range("J115").Value = 40000
$ws.Range("L115").Value = 40000
$ws.Range("N115").Value = -43134
$ws.Range("H122").Value = 3254.7334
$ws.Range("J122").Value = 2024.5
$ws.Range("L122").Value = 6073.5
$ws.Range("N122").Value = -10973.5
$ws.Range("H136").Value = 2989
$ws.Range("I136").Value = 2843.6667
$ws.Range("K136").Value = 8531.000100000001
$ws.Range("M136").Value = -5981.000100000001

Write-Output "Applied 201 cell updates across 7 sheets."
